$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 0.99999999571140963
$ws.Range("A2").Value = 0.99847826396817507
$ws.Range("A3").Value = 0.99454210293986511
$ws.Range("A4").Value = 0.99749455675504506
$ws.Range("A5").Value = 0.98729255514584791
$ws.Range("A6").Value = 0.96480217850274241
$ws.Range("A7").Value = 0.96153736771128462
$ws.Range("A8").Value = 0.95654295530792055
$ws.Range("A9").Value = 0.94479836781111737
$ws.Range("A10").Value = 0.93397042130395613
$ws.Range("A11").Value = 0.93237431155345429
$ws.Range("A12").Value = 0.92959454046893253
$ws.Range("A13").Value = 0.91830761508414094
$ws.Range("A14").Value = 0.91414024835183438
$ws.Range("A15").Value = 0.91154874133884511
$ws.Range("A16").Value = 0.90904223699119868
$ws.Range("A17").Value = 0.9053342813666494
$ws.Range("A18").Value = 0.90422535677000582
$ws.Range("A19").Value = 0.99730323367622553
$ws.Range("A20").Value = 0.9889827947935842
$ws.Range("A21").Value = 0.98758429122786118
$ws.Range("A22").Value = 0.98631978307968571
$ws.Range("A23").Value = 0.98750796000202135
$ws.Range("A24").Value = 0.97448775575893265
$ws.Range("A25").Value = 0.96803086555200557
$ws.Range("A26").Value = 0.96193411531836204
$ws.Range("A27").Value = 0.95818404770645138
$ws.Range("A28").Value = 0.94456261010248932
$ws.Range("A29").Value = 0.92929637569718615
$ws.Range("A30").Value = 0.9227269941320746
$ws.Range("A31").Value = 0.91507387824273723
$ws.Range("A32").Value = 0.9133945571808253
$ws.Range("A33").Value = 0.91287454139558055
